# Applies the "May 2021" carry-forward update:
#  - "April 2021-22" sheet: a couple of formula tweaks + newly-entered
#    "Inter Trf." (column C) figures for a handful of items.
#  - "May 2021" sheet: its external-workbook references move from the
#    (still-empty) "May 2021" tab to the now-populated "April 2021-22" tab,
#    and the "Opening" column (B) is hand-updated to the April closing
#    balances ("Book Balance", column K).

$wb = $excel.ActiveWorkbook

$april = $wb.Worksheets.Item("April 2021-22")
$may   = $wb.Worksheets.Item("May 2021")

# ---------------------------------------------------------------------------
# "April 2021-22" sheet
# ---------------------------------------------------------------------------

# SPECIAL 250GM LAMINATES - wastage % revised 2.3% -> 2.4%
$april.Range("G4").Formula = "=F4*0.024"

# SPECIAL 500 GM POUCH ("RS 5/- LAMINATES " row) - Inter Trf. entered,
# and the wastage figure now derives from a percentage instead of a
# hardcoded sum.
$april.Range("C6").Value = 1086.95
$april.Range("G6").Formula = "=F6*0.108"

# RS 10/-LAMINATES row - wastage figure now derives from a percentage
# instead of a hardcoded sum.
$april.Range("G7").Formula = "=F7*0.075"

# SPECIAL 500 GM POUCH row - Inter Trf. entered.
$april.Range("C11").Value = 65700

# 5 KG PP NON WOVEN Bags row - Inter Trf. entered.
$april.Range("C13").Formula = "=1800+2700"

# 250 GM SPECEIAL CANVES BAGS (24kg) row - Inter Trf. entered.
$april.Range("C25").Value = 2700

# 500GM MUKTA JAR row - Inter Trf. entered.
$april.Range("C29").Value = 32400

# ---------------------------------------------------------------------------
# "May 2021" sheet
# ---------------------------------------------------------------------------

# All the cross-workbook pulls that used to read from the (blank) "May 2021"
# tab of the external workbook now read from "April 2021-22" instead.
$mayExternalRefs = @{
    "F4"  = "='[1]April 2021-22'!`$J`$32"
    "F5"  = "='[1]April 2021-22'!`$C`$32"
    "F6"  = "='[1]April 2021-22'!`$G`$32"
    "F7"  = "='[1]April 2021-22'!`$F`$32"
    "F8"  = "='[1]April 2021-22'!`$E`$32"
    "F10" = "='[1]April 2021-22'!`$H`$32"
    "F11" = "='[1]April 2021-22'!`$I`$32"
    "F12" = "='[1]April 2021-22'!`$B`$32"
    "F13" = "='[1]April 2021-22'!`$D`$32"
    "F29" = "='[1]April 2021-22'!`$K`$32"
    "F30" = "='[1]April 2021-22'!`$L`$32"
}
foreach ($addr in $mayExternalRefs.Keys) {
    $may.Range($addr).Formula = $mayExternalRefs[$addr]
}

# Same two wastage-formula tweaks as on the April sheet.
$may.Range("G6").Formula = "=F6*0.108"
$may.Range("G7").Formula = "=F7*0.075"

# New Inter Trf. (column C) entries.
$may.Range("C4").Value = 4045270
$may.Range("C13").Value = 2700

# Opening balances (column B) rolled forward from April's closing
# "Book Balance" (column K), each re-typed by hand.
$mayOpening = @{
    "B4"  = -5.9
    "B5"  = 448.22
    "B6"  = 1086.95
    "B7"  = 171.85
    "B8"  = 675.76
    "B10" = 38115
    "B11" = 84021
    "B12" = 8616
    "B13" = 4098
    "B15" = 101880
    "B16" = 9188
    "B17" = 56750
    "B18" = 186964
    "B19" = 22434
    "B21" = 599
    "B22" = 1467
    "B23" = 6266
    "B24" = 906
    "B25" = 4130
    "B26" = 505
    "B27" = 361
    "B29" = 30438
    "B30" = 41420
    "B31" = 1130
    "B32" = 526
    "B33" = 82497
    "B34" = 72124
    "B35" = 34902
    "B36" = 6440
}
foreach ($addr in $mayOpening.Keys) {
    $may.Range($addr).Value = $mayOpening[$addr]
}

$excel.Calculate()
